$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2024-12-19 00:22:21", -0.1216357083804905, -0.001645023349513995, 0.008003743216823284),
    @("2024-12-19 00:22:22", -0.1225596532935742, -0.001768774007383996, 0.008671213163986736),
    @("2024-12-19 00:22:23", -0.1225921865651617, -0.001903052328581995, 0.009331973844351572),
    @("2024-12-19 00:22:24", -0.1232038120710058, -0.001859575103587996, 0.009164269663775065),
    @("2024-12-19 00:22:25", -0.1235258914597216, -0.001935394525055996, 0.009562853361352257)
)

$startRow = 135
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}
